$wb = $excel.ActiveWorkbook

$wsTables = $wb.Worksheets.Item("Table_Names")
$wsTables.Range("A1").Value = "T2"
$wsTables.Range("A2").ClearContents()
$wsTables.Range("A3").ClearContents()

$wsFields = $wb.Worksheets.Item("Field_Names")
$wsFields.Range("A1:A5").ClearContents()
